# Generate Report for handback
# Updates the zh-cn and de-de localization-status sheets to reflect that the
# handoff packages have now been handed back and are in sync with en-US.

$wb = $excel.ActiveWorkbook

function Update-HandbackSheet {
    param($SheetName, $HandbackDateTime2, $HandbackDateTime3)

    $ws = $wb.Worksheets.Item($SheetName)
    $links = @($ws.Hyperlinks)

    # --- Row 2 -----------------------------------------------------------
    $ws.Range("B2").Value = "Handed back: in sync with en-US"

    $srcDisplay2 = $ws.Range("A2").Value()
    $srcUrl2 = $links[0].Address
    $handoffDisplay2 = $ws.Range("C2").Value()
    $handoffUrl2 = $links[1].Address

    $ws.Range("E2").Value = $srcDisplay2
    $ws.Hyperlinks.Add($ws.Range("E2"), $srcUrl2, "", "", $srcDisplay2) | Out-Null

    $ws.Range("F2").Value = $handoffDisplay2
    $ws.Hyperlinks.Add($ws.Range("F2"), $handoffUrl2, "", "", $handoffDisplay2) | Out-Null

    $ws.Range("G2").Value = $HandbackDateTime2

    # --- Row 3 -----------------------------------------------------------
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    $srcDisplay3 = $ws.Range("A3").Value()
    $srcUrl3 = $links[2].Address
    $handoffDisplay3 = $ws.Range("C3").Value()
    $handoffUrl3 = $links[3].Address

    $ws.Range("E3").Value = $srcDisplay3
    $ws.Hyperlinks.Add($ws.Range("E3"), $srcUrl3, "", "", $srcDisplay3) | Out-Null

    $ws.Range("F3").Value = $handoffDisplay3
    $ws.Hyperlinks.Add($ws.Range("F3"), $handoffUrl3, "", "", $handoffDisplay3) | Out-Null

    $ws.Range("G3").Value = $HandbackDateTime3
}

Update-HandbackSheet "zh-cn" "2016-01-26 04:47:43" "2016-01-26 04:47:43"
Update-HandbackSheet "de-de" "2016-01-26 04:48:00" "2016-01-26 04:48:00"
